# Apply updated crypto price/volume figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '70.751.29'
$ws.Range("E2").Value = '  -0.65%  '
# Row 3
$ws.Range("D3").Value = '3.804.62'
$ws.Range("E3").Value = '  -1.12%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  -0.07%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '708.27'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +1.71%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.50'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -1.70%  '
# Row 7
$ws.Range("D7").Value = '3.805.40'
$ws.Range("E7").Value = '  -1.04%  '
# Row 8
$ws.Range("E8").Value = '  +0.04%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.521'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -1.05%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -1.69%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.39'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +1.28%  '
# Row 12
$ws.Range("E12").Value = '  -1.15%  '
# Row 13
$ws.Range("E13").Value = '  -2.02%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.11'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -1.16%  '
# Row 15
$ws.Range("D15").Value = '4.441.64'
$ws.Range("E15").Value = '  -1.20%  '
# Row 16
$ws.Range("D16").Value = '3.896.55'
$ws.Range("E16").Value = '  +1.28%  '
# Row 17
$ws.Range("D17").Value = '70.771.81'
$ws.Range("E17").Value = '  -0.67%  '
# Row 18
$ws.Range("E18").Value = '  -0.06%  '
# Row 19
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.40'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -1.98%  '
# Row 20
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.12'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -1.74%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '495.95'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +0.29%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.61'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -5.19%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.730'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +0.89%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.53'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -0.58%  '
# Row 25
$ws.Range("E25").Value = '  -1.28%  '
# Row 26
$ws.Range("E26").Value = '  -1.92%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.45'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -1.25%  '
# Row 28
$ws.Range("D28").Value = '3.955.18'
$ws.Range("E28").Value = '  -1.19%  '
# Row 29
$ws.Range("E29").Value = '  +0.06%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.04'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -4.98%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.08'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -3.27%  '
# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.23'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -2.64%  '
# Row 33
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.30'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -4.42%  '
# Row 34
$ws.Range("E34").Value = '  -2.24%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.173'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  -3.61%  '
# Row 36
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.10'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -2.20%  '
# Row 37
$ws.Range("B37").Value = 'RenzoRestakedETH'
$ws.Range("C37").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D37").Value = '3.774.39'
$ws.Range("E37").Value = '  -0.66%  '
# Row 38
$ws.Range("B38").Value = 'Binance-PegBSC-USD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -0.06%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.101'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -3.54%  '
# Row 40
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.04'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +1.43%  '
# Row 41
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.31'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -2.52%  '
# Row 42
$ws.Range("E42").Value = '  -1.89%  '
# Row 43
$ws.Range("E43").Value = '  -4.19%  '
# Row 44
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -0.03%  '
# Row 45
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("E45").Value = '  +0.07%  '
# Row 46
$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.000322'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +5.58%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '164.98'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +1.00%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '426.68'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +1.61%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '48.74'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +0.18%  '
# Row 50
$ws.Range("E50").Value = '  -0.79%  '
# Row 51
$ws.Range("E51").Value = '  -1.62%  '
